# "alteracao de comportamento e atualizacoes pontuais"
# Adds the new daily quote row (12/02/2026) to the "Cotacoes" sheet and
# keeps the TJLP/SELIC/CDI columns (L:N) formatted the same way row 7 is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: make sure L7/M7/N7 carry the expected number formats -----------
$ws.Range("L7:M7").NumberFormat = "0.0000%"
$ws.Range("N7").NumberFormat = "0.0000000000"

# --- Row 8: new quotation row ------------------------------------------------
$ws.Range("A8").Value = 46065              # 12/02/2026
$ws.Range("B8").Value = 5.1836
$ws.Range("C8").Value = 5.1856
$ws.Range("F8").Value = 5.2147
$ws.Range("G8").Value = 5.3947
$ws.Range("L8").Value = 0.0919
$ws.Range("M8").Value = 0.15
$ws.Range("N8").Value = 0.0551310642
$ws.Range("O8").Value = "OK 12/02/2026 07:55:30"

# Match the number formats already used on the sheet for these columns.
$ws.Range("A8").NumberFormat = "dd/mm/yyyy"
$ws.Range("B8:C8").NumberFormat = "0.0000"
$ws.Range("F8:G8").NumberFormat = "0.0000"
$ws.Range("L8:M8").NumberFormat = "0.0000%"
$ws.Range("N8").NumberFormat = "0.0000000000"
